$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 85 with the latest blue tank titration / CRM accuracy data
# (CRM opened 2022-05-18), continuing the "% off" shared formula pattern
# used by the rows above.
$ws.Range("A85").Value = 20220518
$ws.Range("B85").Value = 2201.1546499999999
$ws.Range("C85").Value = 2224.4699999999998
$ws.Range("D85").Formula = "=100*(B85-C85)/C85"
$ws.Range("E85").Value = 180
$ws.Range("F85").Value = "CRM OPENED 20220518"

# Move the active selection down to the next empty row, as it would be
# after entering the new row of data.
[void]$ws.Range("A86").Select()
